$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A203").Value = 'GO INVESTOR ALERT: Bronstein, Gewirtz and Grossman, LLC Announces that Grocery Outlet Holding Corp. Investors Have Opportunity to Lead Class Action Lawsuit!'
$ws.Range("A204").Value = 'County sues Dayton’s Water Department for alleged ‘breach of contract’'
$ws.Range("A205").Value = 'Montgomery County sues City of Dayton Water Department for breach of contract'
$ws.Range("A206").Value = 'Republicans Target $100 Billion in Alleged COVID-Era Unemployment Fraud'
$ws.Range("A207").Value = 'Tarver Elementary teacher placed on leave after allegations of misconduct'
$ws.Range("A208").Value = 'Alberta cabinet minister calls for removal of health minister amid corruption scandal'
$ws.Range("A209").Value = 'CEO charged by NCD Police'
$ws.Range("A210").Value = 'AG Jackley States Lawsuit Will Not Eliminate Section 504 School Disability Accommodations'
$ws.Range("A211").Value = 'Wayne County deputy awarded $1.7M in lawsuit against driver who hit hit while directing traffic'
$ws.Range("A212").Value = 'Michigan Attorney General files brief in support of lawsuit to block transgender military ban'
$ws.Range("A213").Value = 'Prof. Irina Manta Discusses Lawsuit Against AI Company for Alleged Copyright Infringement'
$ws.Range("A214").Value = 'Seventh attorney resigns after refusing to dismiss case against NYC Mayor Eric Adams'
$ws.Range("A215").Value = 'Farsley Celtic chairman resigns after abuse from supporters'
$ws.Range("A216").Value = 'Local lawyer arrested for harassment'
